# Replace the Chinese placeholder text used in the slide layouts (and the
# slide master) with plain English placeholder words so that the deck uses
# consistent (Latin) fonts. This mirrors the upstream LibreOffice commit
# "replace Japanese text to get consistent fonts used".

$p = $ppt.ActivePresentation

# NOTE: $p.Slides.Item(1).Master.CustomLayouts.Item(n) does not reliably
# index the individual layouts in this host -- go through
# $p.SlideMaster.CustomLayouts.Item(n) instead, which does.
$master = $p.SlideMaster

function Set-FiveLevelText($shape, $l1, $l2, $l3, $l4, $l5) {
    $tr = $shape.TextFrame.TextRange
    $tr.Paragraphs(1, 1).Text = $l1
    $tr.Paragraphs(2, 1).Text = $l2
    $tr.Paragraphs(3, 1).Text = $l3
    $tr.Paragraphs(4, 1).Text = $l4
    $tr.Paragraphs(5, 1).Text = $l5
}

# --- slideLayout1.xml (title) ---------------------------------------------
$cl = $master.CustomLayouts.Item(1)
$cl.Name = "THING"
$cl.Shapes.Item(1).TextFrame.TextRange.Text = "HELLO"
$cl.Shapes.Item(2).TextFrame.TextRange.Text = "WORLD"

# --- slideLayout2.xml (obj) ------------------------------------------------
$cl = $master.CustomLayouts.Item(2)
$cl.Name = "STYLE"
$cl.Shapes.Item(1).TextFrame.TextRange.Text = "HELLO"
Set-FiveLevelText $cl.Shapes.Item(2) "SOMETHING" "FOO" "BAR" "BAR" "BAZ"

# --- slideLayout3.xml (secHead) --------------------------------------------
$cl = $master.CustomLayouts.Item(3)
$cl.Name = "ONE"
$cl.Shapes.Item(1).TextFrame.TextRange.Text = "HELLO"
$cl.Shapes.Item(2).TextFrame.TextRange.Text = "SOMETHING"

# --- slideLayout4.xml (twoObj) ---------------------------------------------
$cl = $master.CustomLayouts.Item(4)
$cl.Name = "FOUR"
$cl.Shapes.Item(1).TextFrame.TextRange.Text = "HELLO"
Set-FiveLevelText $cl.Shapes.Item(2) "SOMETHING" "FOO" "BAR" "BAR" "BAZ"
Set-FiveLevelText $cl.Shapes.Item(3) "SOMETHING" "FOO" "BAR" "BAR" "BAZ"

# --- slideLayout5.xml (twoTxTwoObj) -----------------------------------------
$cl = $master.CustomLayouts.Item(5)
$cl.Name = "EG"
$cl.Shapes.Item(1).TextFrame.TextRange.Text = "HELLO"
$cl.Shapes.Item(2).TextFrame.TextRange.Text = "SOMETHING"
Set-FiveLevelText $cl.Shapes.Item(3) "SOMETHING" "FOO" "BAR" "BAR" "BAZ"
$cl.Shapes.Item(4).TextFrame.TextRange.Text = "SOMETHING"
Set-FiveLevelText $cl.Shapes.Item(5) "SOMETHING" "FOO" "BAR" "BAR" "BAZ"

# --- slideLayout6.xml (titleOnly) -------------------------------------------
$cl = $master.CustomLayouts.Item(6)
$cl.Name = "ABC"
$cl.Shapes.Item(1).TextFrame.TextRange.Text = "HELLO"

# --- slideLayout7.xml (blank) -----------------------------------------------
$cl = $master.CustomLayouts.Item(7)
$cl.Name = "XY"

# --- slideLayout8.xml (objTx) -----------------------------------------------
$cl = $master.CustomLayouts.Item(8)
$cl.Name = "ABCDEF"
$cl.Shapes.Item(1).TextFrame.TextRange.Text = "HELLO"
Set-FiveLevelText $cl.Shapes.Item(2) "SOMETHING" "FOO" "BAR" "BAR" "BAZ"
$cl.Shapes.Item(3).TextFrame.TextRange.Text = "SOMETHING"

# --- slideLayout9.xml (picTx) -----------------------------------------------
$cl = $master.CustomLayouts.Item(9)
$cl.Name = "HIJKL"
$cl.Shapes.Item(1).TextFrame.TextRange.Text = "HELLO"
$cl.Shapes.Item(3).TextFrame.TextRange.Text = "SOMETHING"

# --- slideLayout10.xml (vertTx) ---------------------------------------------
$cl = $master.CustomLayouts.Item(10)
$cl.Name = "EXAMPLE"
$cl.Shapes.Item(1).TextFrame.TextRange.Text = "HELLO"
Set-FiveLevelText $cl.Shapes.Item(2) "SOMETHING" "FOO" "BAR" "BAR" "BAZ"

# --- slideLayout11.xml (vertTitleAndTx) --------------------------------------
$cl = $master.CustomLayouts.Item(11)
$cl.Name = "SOMETHING"
$cl.Shapes.Item(1).TextFrame.TextRange.Text = "HELLO"
Set-FiveLevelText $cl.Shapes.Item(2) "SOMETHING" "FOO" "BAR" "BAR" "BAZ"

# --- slideMaster1.xml ---------------------------------------------------------
$master.Shapes.Item(1).TextFrame.TextRange.Text = "HELLO"
Set-FiveLevelText $master.Shapes.Item(2) "SOMETHING" "FOO" "BAR" "BAR" "BAZ"
